# Fruta / hortaliza, semanal
# Insert a new weekly record as row 68 on the active sheet, shifting the
# existing rows 68-108 down to 69-109 (dimension grows from A1:T108 to A1:T109).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 68.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row with the new weekly price entry.
$ws.Cells.Item(68, 1).Value  = 4
$ws.Cells.Item(68, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(68, 3).Value  = "Los Lagos"
$ws.Cells.Item(68, 4).Value  = 44813
$ws.Cells.Item(68, 5).Value  = 10
$ws.Cells.Item(68, 6).Value  = "Fruta"
$ws.Cells.Item(68, 7).Value  = 100104
$ws.Cells.Item(68, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(68, 9).Value  = 100104003
$ws.Cells.Item(68, 10).Value = "Membrillo"
$ws.Cells.Item(68, 11).Value = "Champion"
$ws.Cells.Item(68, 12).Value = "Primera"
$ws.Cells.Item(68, 13).Value = 300
$ws.Cells.Item(68, 14).Value = 14000
$ws.Cells.Item(68, 15).Value = 15000
$ws.Cells.Item(68, 16).Value = 14500
$ws.Cells.Item(68, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(68, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(68, 19).Value = 806
$ws.Cells.Item(68, 20).Value = 18
